$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '40.074.17'
$ws.Range("E2").Value = '  -0.07%  '
$ws.Range("D3").Value = '2.212.05'
$ws.Range("E3").Value = '  -1.18%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '294.85'
$ws.Range("E5").Value = '  -0.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '88.25'
$ws.Range("E6").Value = '  +1.56%  '
$ws.Range("E7").Value = '  -1.14%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.475'
$ws.Range("E9").Value = '  +0.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '30.17'
$ws.Range("E10").Value = '  -3.77%  '
$ws.Range("E11").Value = '  -1.71%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '49.87'
$ws.Range("E12").Value = '  +5.87%  '
$ws.Range("E13").Value = '  +2.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.48'
$ws.Range("E14").Value = '  +0.38%  '
$ws.Range("D15").Value = '2.553.74'
$ws.Range("E15").Value = '  -0.95%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.77'
$ws.Range("E16").Value = '  -2.90%  '
$ws.Range("D17").Value = '2.176.30'
$ws.Range("E17").Value = '  -4.08%  '
$ws.Range("E18").Value = '  -0.89%  '
$ws.Range("D19").Value = '39.972.92'
$ws.Range("E19").Value = '  -0.13%  '
$ws.Range("D20").Value = '0.0₃0885'
$ws.Range("E20").Value = '  -0.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.35'
$ws.Range("E21").Value = '  +3.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.78'
$ws.Range("E22").Value = '  -0.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.29'
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '237.60'
$ws.Range("E24").Value = '  +0.81%  '
$ws.Range("E25").Value = '  +0.12%  '
$ws.Range("E26").Value = '  -0.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.81'
$ws.Range("E27").Value = '  -1.97%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.50'
$ws.Range("E28").Value = '  -1.83%  '
$ws.Range("E29").Value = '  -3.08%  '
$ws.Range("E30").Value = '  -0.96%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '156.66'
$ws.Range("E31").Value = '  +2.76%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.45'
$ws.Range("E32").Value = '  -5.85%  '
$ws.Range("E33").Value = '  +0.08%  '
$ws.Range("E34").Value = '  +0.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0712'
$ws.Range("E35").Value = '  -0.75%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.34'
$ws.Range("E37").Value = '  +2.97%  '
$ws.Range("E38").Value = '  +0.51%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0980'
$ws.Range("E39").Value = '  -2.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '15.43'
$ws.Range("E40").Value = '  -5.56%  '
$ws.Range("E41").Value = '  -2.09%  '
$ws.Range("D42").Value = '2.122.17'
$ws.Range("E42").Value = '  +3.79%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.75'
$ws.Range("E43").Value = '  -2.34%  '
$ws.Range("E44").Value = '  -2.47%  '
$ws.Range("E45").Value = '  -0.97%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.45'
$ws.Range("E46").Value = '  +5.86%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.64'
$ws.Range("E47").Value = '  -4.57%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.66'
$ws.Range("E48").Value = '  +3.32%  '
$ws.Range("D49").Value = '2.422.19'
$ws.Range("E49").Value = '  -1.06%  '
$ws.Range("E50").Value = '  +2.48%  '
$ws.Range("E51").Value = '  +0.88%  '
